# "refine figures and ch1"
# Shrinks a handful of flow-diagram text boxes / the connector between them
# on the single slide of ch3_TestSelection.pptx: smaller font (14pt -> 12pt)
# and tighter box sizes, plus a couple of position tweaks.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 : "大规模无标注测试数据"  (textbox id=39) ---------------------
$sh1 = $s.Shapes.Item(1)
$sh1.TextFrame.TextRange.Font.Size = 12
$sh1.Height = 36.25

# --- Shape 2 : "决策路径"  (textbox id=63) ---------------------------------
$sh2 = $s.Shapes.Item(2)
$sh2.TextFrame.TextRange.Font.Size = 12
$sh2.Left = 235.9000016
$sh2.Width = 62.4
$sh2.Height = 21.7

# --- Shape 7 : "代表性数据选取"  (textbox id=93) ---------------------------
$sh7 = $s.Shapes.Item(7)
$sh7.TextFrame.TextRange.Font.Size = 12
$sh7.Height = 21.7

# --- Shape 9 : bent connector between the boxes (id=264) -------------------
$sh9 = $s.Shapes.Item(9)
$sh9.Left = 483.8000031
$sh9.Top = 305.3000031
$sh9.Width = 19.35
$sh9.Height = 139.85

# --- Shape 20 : "边界数据选取"  (textbox id=174) ----------------------------
$sh20 = $s.Shapes.Item(20)
$sh20.TextFrame.TextRange.Font.Size = 12
$sh20.Left = 336.5500031
$sh20.Top = 374.0500031
$sh20.Width = 87.0
$sh20.Height = 21.7

# --- Shape 28 : "测试集生成"  (textbox id=57) -------------------------------
$sh28 = $s.Shapes.Item(28)
$sh28.TextFrame.TextRange.Font.Size = 12
$sh28.Height = 21.7
